$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2717632.2
$ws.Range("J28").Value = 7581.3335
$ws.Range("L28").Value = 7581.3335
$ws.Range("N28").Value = -8551.333500000001
$ws.Range("H40").Value = 3000
$ws.Range("I40").Value = 3000
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 3000
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2825
$ws.Range("N40").ClearContents()
$ws.Range("H86").Value = 1624.375
$ws.Range("H89").Value = 1624.375
$ws.Range("H98").Value = 2079.9473
$ws.Range("I98").Value = 2330.0625
$ws.Range("K98").Value = 2330.0625
$ws.Range("M98").Value = -832.0625
$ws.Range("H122").Value = 2079.9473
$ws.Range("I122").Value = 2330.0625
$ws.Range("K122").Value = 6990.1875
$ws.Range("M122").Value = -4540.1875
$ws.Range("H129").Value = 860.8333
$ws.Range("J129").Value = 884.0714
$ws.Range("L129").Value = 2652.2142
$ws.Range("N129").Value = -12652.2142
$ws.Range("H132").Value = 908.4464
$ws.Range("I132").Value = 771.9216
$ws.Range("K132").Value = 2315.7648
$ws.Range("M132").Value = 214.2352000000001
$ws.Range("H137").Value = 37196.465
$ws.Range("I137").Value = 1415.45
$ws.Range("J137").Value = 126649
$ws.Range("K137").Value = 4246.35
$ws.Range("L137").Value = 379947
$ws.Range("M137").Value = -1696.35
$ws.Range("N137").Value = -385047
$ws.Range("H141").Value = 1039809.06
$ws.Range("I141").Value = 1335246.9
$ws.Range("K141").Value = 4005740.7
$ws.Range("M141").Value = -4000560.7

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2940.494
$ws.Range("I32").Value = 2281.2388
$ws.Range("K32").Value = 2281.2388
$ws.Range("M32").Value = -1994.2388
$ws.Range("H45").Value = 1657.7142
$ws.Range("I45").Value = 1137.75
$ws.Range("K45").Value = 1137.75
$ws.Range("M45").Value = -760.75
$ws.Range("H61").Value = 2046.2
$ws.Range("I61").Value = 1096.8334
$ws.Range("J61").Value = 5843.6665
$ws.Range("K61").Value = 1096.8334
$ws.Range("L61").Value = 5843.6665
$ws.Range("M61").Value = -884.8334
$ws.Range("N61").Value = -6267.6665
$ws.Range("H74").Value = 1186.8043
$ws.Range("I74").Value = 999.3333
$ws.Range("J74").Value = 1861.7
$ws.Range("K74").Value = 999.3333
$ws.Range("L74").Value = 1861.7
$ws.Range("M74").Value = -125.3333
$ws.Range("N74").Value = -3609.7
$ws.Range("H77").Value = 1186.8043
$ws.Range("I77").Value = 999.3333
$ws.Range("J77").Value = 1861.7
$ws.Range("K77").Value = 4996.6665
$ws.Range("L77").Value = 9308.5
$ws.Range("M77").Value = -628.6665000000003
$ws.Range("N77").Value = -18044.5
$ws.Range("H97").Value = 804.1
$ws.Range("I97").Value = 620.1429000000001
$ws.Range("J97").Value = 1233.3334
$ws.Range("K97").Value = 620.1429000000001
$ws.Range("L97").Value = 1233.3334
$ws.Range("M97").Value = -124.1429000000001
$ws.Range("N97").Value = -2225.3334
$ws.Range("H102").Value = 1809.0476
$ws.Range("I102").Value = 1528.8823
$ws.Range("K102").Value = 1528.8823
$ws.Range("M102").Value = 93.11770000000001
$ws.Range("H110").Value = 2674.4285
$ws.Range("I110").Value = 1677.75
$ws.Range("K110").Value = 1677.75
$ws.Range("M110").Value = 367.25
$ws.Range("H122").Value = 3140
$ws.Range("I122").Value = 1996.6666
$ws.Range("K122").Value = 5989.9998
$ws.Range("M122").Value = -3539.9998
$ws.Range("H136").Value = 2046.2
$ws.Range("I136").Value = 1096.8334
$ws.Range("J136").Value = 5843.6665
$ws.Range("K136").Value = 3290.5002
$ws.Range("L136").Value = 17530.9995
$ws.Range("M136").Value = -740.5001999999999
$ws.Range("N136").Value = -22630.9995

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 8694.583000000001
$ws.Range("I80").Value = 349
$ws.Range("J80").Value = 11476.444
$ws.Range("K80").Value = 349
$ws.Range("L80").Value = 11476.444
$ws.Range("M80").Value = 649
$ws.Range("N80").Value = -13472.444
$ws.Range("H83").Value = 8694.583000000001
$ws.Range("I83").Value = 349
$ws.Range("J83").Value = 11476.444
$ws.Range("K83").Value = 1745
$ws.Range("L83").Value = 57382.22
$ws.Range("M83").Value = 3247
$ws.Range("N83").Value = -67366.22
$ws.Range("H107").Value = 1828.3684
$ws.Range("I107").Value = 1767.1428
$ws.Range("J107").Value = 1999.8
$ws.Range("K107").Value = 1767.1428
$ws.Range("L107").Value = 1999.8
$ws.Range("M107").Value = 152.8571999999999
$ws.Range("N107").Value = -5839.8

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1143.5
$ws.Range("I22").Value = 400
$ws.Range("J22").Value = 1391.3334
$ws.Range("K22").Value = 400
$ws.Range("L22").Value = 1391.3334
$ws.Range("M22").Value = -50
$ws.Range("N22").Value = -2091.3334
$ws.Range("H31").Value = 2174.4167
$ws.Range("I31").Value = 1785.2858
$ws.Range("J31").Value = 2719.2
$ws.Range("K31").Value = 1785.2858
$ws.Range("L31").Value = 2719.2
$ws.Range("M31").Value = -1490.2858
$ws.Range("N31").Value = -3309.2
$ws.Range("H34").Value = 2174.4167
$ws.Range("I34").Value = 1785.2858
$ws.Range("J34").Value = 2719.2
$ws.Range("K34").Value = 1785.2858
$ws.Range("L34").Value = 2719.2
$ws.Range("M34").Value = -1583.2858
$ws.Range("N34").Value = -3123.2
$ws.Range("H122").Value = 4064.077
$ws.Range("I122").Value = 2517
$ws.Range("K122").Value = 7551
$ws.Range("M122").Value = -5101
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H132").Value = 1741.28
$ws.Range("I132").Value = 1098.6666
$ws.Range("K132").Value = 3295.9998
$ws.Range("M132").Value = -765.9998000000001

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 803.5714
$ws.Range("J5").Value = 805
$ws.Range("L5").Value = 2415
$ws.Range("N5").Value = -2639
$ws.Range("H34").Value = 507.8
$ws.Range("I34").Value = 137.5
$ws.Range("J34").Value = 1989
$ws.Range("K34").Value = 412.5
$ws.Range("L34").Value = 5967
$ws.Range("M34").Value = -328.5
$ws.Range("N34").Value = -6135
$ws.Range("H39").Value = 6000
$ws.Range("J39").Value = 6000
$ws.Range("L39").Value = 18000
$ws.Range("N39").Value = -18588
$ws.Range("H55").Value = 3200
$ws.Range("J55").Value = 3632.6667
$ws.Range("L55").Value = 10898.0001
$ws.Range("N55").Value = -11252.0001
$ws.Range("H131").Value = 814.73
$ws.Range("I131").Value = 800
$ws.Range("J131").Value = 814.8788
$ws.Range("K131").Value = 2400
$ws.Range("L131").Value = 2444.6364
$ws.Range("M131").Value = 2640
$ws.Range("N131").Value = -12524.6364
$ws.Range("H133").Value = 1399
$ws.Range("I133").Value = 1399
$ws.Range("K133").Value = 4197
$ws.Range("M133").Value = 863
$ws.Range("H135").Value = 803.5714
$ws.Range("J135").Value = 805
$ws.Range("L135").Value = 7245
$ws.Range("N135").Value = -12315
$ws.Range("H139").Value = 21800.6
$ws.Range("I139").Value = 100000
$ws.Range("K139").Value = 300000
$ws.Range("M139").Value = -294860

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1284649.5
$ws.Range("I132").Value = 1925081
$ws.Range("K132").Value = 5775243
$ws.Range("M132").Value = -5772713

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2757.625
$ws.Range("I7").Value = 1769.4166
$ws.Range("J7").Value = 5722.25
$ws.Range("K7").Value = 1769.4166
$ws.Range("L7").Value = 5722.25
$ws.Range("M7").Value = -1657.4166
$ws.Range("N7").Value = -5946.25
$ws.Range("H126").Value = 2757.625
$ws.Range("I126").Value = 1769.4166
$ws.Range("J126").Value = 5722.25
$ws.Range("K126").Value = 5308.2498
$ws.Range("L126").Value = 17166.75
$ws.Range("M126").Value = -2838.2498
$ws.Range("N126").Value = -22106.75
$ws.Range("H136").Value = 2307.5356
$ws.Range("I136").Value = 1430.55
$ws.Range("J136").Value = 4500
$ws.Range("K136").Value = 4291.65
$ws.Range("L136").Value = 13500
$ws.Range("M136").Value = -1741.65
$ws.Range("N136").Value = -18600

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1200
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 1200
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 1200
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -3946
$ws.Range("H122").Value = 61453.77
$ws.Range("I122").Value = 112167.86
$ws.Range("K122").Value = 336503.58
$ws.Range("M122").Value = -334053.58
$ws.Range("H132").Value = 1172.8684
$ws.Range("I132").Value = 850.5172
$ws.Range("J132").Value = 2211.5557
$ws.Range("K132").Value = 2551.5516
$ws.Range("L132").Value = 6634.6671
$ws.Range("M132").Value = -21.55159999999978
$ws.Range("N132").Value = -11694.6671
